# Tsalka.xlsx - "upgrade left table until javakheti"
#
# The left-hand statistics table (rows 5-7: Total / Urban / Rural) gets
# several of its yearly figures replaced with the "confidential / not
# available" placeholder "..." (a literal three-dot string, distinct from
# the single-character ellipsis "…" already used elsewhere in the sheet),
# the sheet gets its real name "Tsalka" instead of the generic "1", and
# the stray blank row between the table and the footnote is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Give the worksheet its proper name.
$ws.Name = "Tsalka"

# 2. Row 6 ("Urban"): every year from 2011 to 2021 becomes confidential,
#    as well as 2023. 2018 and 2022 already show the "…" placeholder and
#    stay untouched.
$ws.Range("C6:M6").Value = "..."
$ws.Range("O6").Value = "..."

# 3. Row 7 ("Rural"): 2012-2014 become confidential; 2021 and 2023 switch
#    from real numbers to the placeholder values ("…" and "..." resp.).
$ws.Range("D7:F7").Value = "..."
$ws.Range("M7").Value = "…"
$ws.Range("O7").Value = "..."

# 4. Remove the empty row that used to sit between the table (row 7) and
#    the footnote (was row 9), so the footnote becomes row 8.
$ws.Rows(8).Delete()
